$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 152
$ws.Range("I9").Value = 121.92857
$ws.Range("J9").Value = 198.77777
$ws.Range("K9").Value = 121.92857
$ws.Range("L9").Value = 198.77777
$ws.Range("M9").Value = 47.07143000000001
$ws.Range("N9").Value = -536.77777
$ws.Range("H33").Value = 4683.9565
$ws.Range("I33").Value = 5847.278
$ws.Range("J33").Value = 496
$ws.Range("K33").Value = 5847.278
$ws.Range("L33").Value = 496
$ws.Range("M33").Value = -5618.278
$ws.Range("N33").Value = -954
$ws.Range("H62").Value = 3556
$ws.Range("I62").Value = 3072.6667
$ws.Range("K62").Value = 3072.6667
$ws.Range("M62").Value = -2448.6667
$ws.Range("H65").Value = 3556
$ws.Range("I65").Value = 3072.6667
$ws.Range("K65").Value = 15363.3335
$ws.Range("M65").Value = -12243.3335
$ws.Range("H87").Value = 18175.598
$ws.Range("J87").Value = 18175.598
$ws.Range("L87").Value = 18175.598
$ws.Range("N87").Value = -20671.598
$ws.Range("H90").Value = 18175.598
$ws.Range("J90").Value = 18175.598
$ws.Range("L90").Value = 54526.79400000001
$ws.Range("N90").Value = -67006.79400000001
$ws.Range("H107").Value = 561641.6
$ws.Range("I107").Value = 721688.9399999999
$ws.Range("J107").Value = 1476
$ws.Range("K107").Value = 721688.9399999999
$ws.Range("L107").Value = 1476
$ws.Range("M107").Value = -719768.9399999999
$ws.Range("N107").Value = -5316
$ws.Range("H129").Value = 1260.963
$ws.Range("I129").Value = 835
$ws.Range("J129").Value = 1295.04
$ws.Range("K129").Value = 2505
$ws.Range("L129").Value = 3885.12
$ws.Range("M129").Value = 2495
$ws.Range("N129").Value = -13885.12
$ws.Range("H140").Value = 60782.637
$ws.Range("J140").Value = 60782.637
$ws.Range("L140").Value = 60782.637
$ws.Range("N140").Value = -71142.637

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1628.0526
$ws.Range("I2").Value = 1874.5333
$ws.Range("J2").Value = 703.75
$ws.Range("K2").Value = 1874.5333
$ws.Range("L2").Value = 703.75
$ws.Range("M2").Value = -1761.5333
$ws.Range("N2").Value = -929.75
$ws.Range("H32").Value = 4767.511
$ws.Range("I32").Value = 3372.9768
$ws.Range("K32").Value = 3372.9768
$ws.Range("M32").Value = -3085.9768
$ws.Range("H110").Value = 26759.6
$ws.Range("I110").Value = 34680.8
$ws.Range("J110").Value = 2996
$ws.Range("K110").Value = 34680.8
$ws.Range("L110").Value = 2996
$ws.Range("M110").Value = -32635.8
$ws.Range("N110").Value = -7086
$ws.Range("H116").Value = 1628.0526
$ws.Range("I116").Value = 1874.5333
$ws.Range("J116").Value = 703.75
$ws.Range("K116").Value = 1874.5333
$ws.Range("L116").Value = 703.75
$ws.Range("M116").Value = 419.4666999999999
$ws.Range("N116").Value = -5291.75

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1628.0526
$ws.Range("I3").Value = 1874.5333
$ws.Range("J3").Value = 703.75
$ws.Range("K3").Value = 1874.5333
$ws.Range("L3").Value = 703.75
$ws.Range("M3").Value = -1760.5333
$ws.Range("N3").Value = -931.75
$ws.Range("H107").Value = 836
$ws.Range("I107").Value = 680
$ws.Range("J107").Value = 1148
$ws.Range("K107").Value = 680
$ws.Range("L107").Value = 1148
$ws.Range("M107").Value = 1240
$ws.Range("N107").Value = -4988

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2566.923
$ws.Range("I16").Value = 2677.9092
$ws.Range("K16").Value = 2677.9092
$ws.Range("M16").Value = -2390.9092
$ws.Range("H107").Value = 868.4545000000001
$ws.Range("I107").Value = 583.36365
$ws.Range("J107").Value = 1153.5454
$ws.Range("K107").Value = 583.36365
$ws.Range("L107").Value = 1153.5454
$ws.Range("M107").Value = 1336.63635
$ws.Range("N107").Value = -4993.5454
$ws.Range("H113").Value = 2566.923
$ws.Range("I113").Value = 2677.9092
$ws.Range("K113").Value = 2677.9092
$ws.Range("M113").Value = -507.9092000000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 818.5
$ws.Range("I114").Value = 336
$ws.Range("J114").Value = 1163.1428
$ws.Range("K114").Value = 1008
$ws.Range("L114").Value = 3489.4284
$ws.Range("M114").Value = 2246
$ws.Range("N114").Value = -9997.428400000001
$ws.Range("H137").Value = 2115.4688
$ws.Range("I137").Value = 2403.2222
$ws.Range("J137").Value = 2002.8695
$ws.Range("K137").Value = 7209.6666
$ws.Range("L137").Value = 6008.6085
$ws.Range("M137").Value = -2109.6666
$ws.Range("N137").Value = -16208.6085

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H107").Value = 1272.4445
$ws.Range("I107").Value = 1593.1428
$ws.Range("J107").Value = 150
$ws.Range("K107").Value = 1593.1428
$ws.Range("L107").Value = 150
$ws.Range("M107").Value = 326.8571999999999
$ws.Range("N107").Value = -3990
$ws.Range("H113").Value = 2251.9
$ws.Range("I113").Value = 2064.875
$ws.Range("K113").Value = 2064.875
$ws.Range("M113").Value = 105.125

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1428
$ws.Range("I68").Value = 590
$ws.Range("J68").Value = 1986.6666
$ws.Range("K68").Value = 590
$ws.Range("L68").Value = 1986.6666
$ws.Range("M68").Value = 159
$ws.Range("N68").Value = -3484.6666
$ws.Range("H71").Value = 1428
$ws.Range("I71").Value = 590
$ws.Range("J71").Value = 1986.6666
$ws.Range("K71").Value = 2950
$ws.Range("L71").Value = 9933.333000000001
$ws.Range("M71").Value = 794
$ws.Range("N71").Value = -17421.333
$ws.Range("H136").Value = 10102714
$ws.Range("I136").Value = 1800.9
$ws.Range("J136").Value = 111111850
$ws.Range("K136").Value = 5402.700000000001
$ws.Range("L136").Value = 333335550
$ws.Range("M136").Value = -2852.700000000001
$ws.Range("N136").Value = -333340650

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H132").Value = 2007.1936
$ws.Range("I132").Value = 2081.6538
$ws.Range("J132").Value = 1620
$ws.Range("K132").Value = 6244.9614
$ws.Range("L132").Value = 4860
$ws.Range("M132").Value = -3714.9614
$ws.Range("N132").Value = -9920
